$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 2537.75
$ws.Cells.Item(5, 9).Value = 2896
$ws.Cells.Item(5, 11).Value = 2896
$ws.Cells.Item(5, 13).Value = -2781
$ws.Cells.Item(17, 8).Value = 4500
$ws.Cells.Item(17, 10).Value = 6000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 14).Value = -18336
$ws.Cells.Item(26, 8).Value = 1199.5
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).Value = $null
$ws.Cells.Item(28, 8).Value = 534.9375
$ws.Cells.Item(28, 10).Value = 1002
$ws.Cells.Item(28, 12).Value = 1002
$ws.Cells.Item(28, 14).Value = -1972
$ws.Cells.Item(41, 8).Value = 615.4
$ws.Cells.Item(41, 9).Value = 569.25
$ws.Cells.Item(41, 10).Value = 800
$ws.Cells.Item(41, 11).Value = 569.25
$ws.Cells.Item(41, 12).Value = 800
$ws.Cells.Item(41, 13).Value = -129.25
$ws.Cells.Item(41, 14).Value = -1680
$ws.Cells.Item(131, 8).Value = 1541.8
$ws.Cells.Item(131, 10).Value = 1337.5
$ws.Cells.Item(131, 12).Value = 4012.5
$ws.Cells.Item(131, 14).Value = -14092.5
$ws.Cells.Item(132, 8).Value = 2814.4614
$ws.Cells.Item(132, 9).Value = 2549.8333
$ws.Cells.Item(132, 11).Value = 7649.499899999999
$ws.Cells.Item(132, 13).Value = -5119.499899999999
$ws.Cells.Item(141, 8).Value = 5805.2144
$ws.Cells.Item(141, 9).Value = 5482.769
$ws.Cells.Item(141, 11).Value = 16448.307
$ws.Cells.Item(141, 13).Value = -11268.307

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 4999.5
$ws.Cells.Item(22, 10).Value = 4999.5
$ws.Cells.Item(22, 12).Value = 4999.5
$ws.Cells.Item(22, 14).Value = -5597.5
$ws.Cells.Item(37, 8).Value = 24497.9
$ws.Cells.Item(37, 9).Value = 20000
$ws.Cells.Item(37, 11).Value = 20000
$ws.Cells.Item(37, 13).Value = -19727
$ws.Cells.Item(74, 8).Value = 4999.5
$ws.Cells.Item(74, 9).Value = 4999.5
$ws.Cells.Item(74, 11).Value = 4999.5
$ws.Cells.Item(74, 13).Value = -4125.5
$ws.Cells.Item(77, 8).Value = 4999.5
$ws.Cells.Item(77, 9).Value = 4999.5
$ws.Cells.Item(77, 11).Value = 25000
$ws.Cells.Item(77, 13).Value = -20629.5
$ws.Cells.Item(132, 8).Value = 2782.75
$ws.Cells.Item(132, 9).Value = 2568
$ws.Cells.Item(132, 11).Value = 7704
$ws.Cells.Item(132, 13).Value = -5174

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1306.3182
$ws.Cells.Item(94, 9).Value = 1333.1052
$ws.Cells.Item(94, 10).Value = 1136.6666
$ws.Cells.Item(94, 11).Value = 1333.1052
$ws.Cells.Item(94, 12).Value = 1136.6666
$ws.Cells.Item(94, 13).Value = -882.1052
$ws.Cells.Item(94, 14).Value = -2038.6666
$ws.Cells.Item(134, 8).Value = 3795
$ws.Cells.Item(134, 9).Value = 3834.5
$ws.Cells.Item(134, 11).Value = 11503.5
$ws.Cells.Item(134, 13).Value = -8968.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 20410.25
$ws.Cells.Item(15, 9).Value = 20410.25
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 20410.25
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -20240.25
$ws.Cells.Item(15, 14).Value = $null
$ws.Cells.Item(31, 8).Value = 4546.1763
$ws.Cells.Item(31, 9).Value = 4205.5454
$ws.Cells.Item(31, 10).Value = 5170.6665
$ws.Cells.Item(31, 11).Value = 4205.5454
$ws.Cells.Item(31, 12).Value = 5170.6665
$ws.Cells.Item(31, 13).Value = -3910.5454
$ws.Cells.Item(31, 14).Value = -5760.6665
$ws.Cells.Item(34, 8).Value = 4546.1763
$ws.Cells.Item(34, 9).Value = 4205.5454
$ws.Cells.Item(34, 10).Value = 5170.6665
$ws.Cells.Item(34, 11).Value = 4205.5454
$ws.Cells.Item(34, 12).Value = 5170.6665
$ws.Cells.Item(34, 13).Value = -4003.5454
$ws.Cells.Item(34, 14).Value = -5574.6665
$ws.Cells.Item(69, 8).Value = 4500
$ws.Cells.Item(69, 9).Value = 4500
$ws.Cells.Item(69, 11).Value = 4500
$ws.Cells.Item(69, 13).Value = -3751
$ws.Cells.Item(72, 8).Value = 4500
$ws.Cells.Item(72, 9).Value = 4500
$ws.Cells.Item(72, 11).Value = 13500
$ws.Cells.Item(72, 13).Value = -9756
$ws.Cells.Item(99, 8).Value = 3547.3333
$ws.Cells.Item(99, 9).Value = 1956.5
$ws.Cells.Item(99, 10).Value = 4001.8572
$ws.Cells.Item(99, 11).Value = 1956.5
$ws.Cells.Item(99, 12).Value = 4001.8572
$ws.Cells.Item(99, 13).Value = -458.5
$ws.Cells.Item(99, 14).Value = -6997.8572
$ws.Cells.Item(126, 8).Value = 3547.3333
$ws.Cells.Item(126, 9).Value = 1956.5
$ws.Cells.Item(126, 10).Value = 4001.8572
$ws.Cells.Item(126, 11).Value = 5869.5
$ws.Cells.Item(126, 12).Value = 12005.5716
$ws.Cells.Item(126, 13).Value = -3399.5
$ws.Cells.Item(126, 14).Value = -16945.5716
$ws.Cells.Item(132, 8).Value = 2394.4666
$ws.Cells.Item(132, 9).Value = 2075.3333
$ws.Cells.Item(132, 10).Value = 3671
$ws.Cells.Item(132, 11).Value = 6225.999899999999
$ws.Cells.Item(132, 12).Value = 11013
$ws.Cells.Item(132, 13).Value = -3695.999899999999
$ws.Cells.Item(132, 14).Value = -16073

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 1103.1428
$ws.Cells.Item(8, 9).Value = 1103.1428
$ws.Cells.Item(8, 11).Value = 3309.4284
$ws.Cells.Item(8, 13).Value = -3170.4284
$ws.Cells.Item(34, 8).Value = 1372.2858
$ws.Cells.Item(34, 9).Value = 301.5
$ws.Cells.Item(34, 10).Value = 2800
$ws.Cells.Item(34, 11).Value = 904.5
$ws.Cells.Item(34, 12).Value = 8400
$ws.Cells.Item(34, 13).Value = -820.5
$ws.Cells.Item(34, 14).Value = -8568
$ws.Cells.Item(121, 8).Value = 895.5
$ws.Cells.Item(121, 9).Value = 100
$ws.Cells.Item(121, 10).Value = 1009.1429
$ws.Cells.Item(121, 11).Value = 300
$ws.Cells.Item(121, 12).Value = 3027.4287
$ws.Cells.Item(121, 13).Value = 1010
$ws.Cells.Item(121, 14).Value = -5647.4287
$ws.Cells.Item(131, 8).Value = 1990.579
$ws.Cells.Item(131, 9).Value = 1202.2
$ws.Cells.Item(131, 10).Value = 2866.5557
$ws.Cells.Item(131, 11).Value = 3606.6
$ws.Cells.Item(131, 12).Value = 8599.667099999999
$ws.Cells.Item(131, 13).Value = 1433.4
$ws.Cells.Item(131, 14).Value = -18679.6671
$ws.Cells.Item(138, 8).Value = 1999
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 13).Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 39000
$ws.Cells.Item(20, 10).Value = 39000
$ws.Cells.Item(20, 12).Value = 39000
$ws.Cells.Item(20, 14).Value = -39490
$ws.Cells.Item(62, 8).Value = 27000
$ws.Cells.Item(62, 9).Value = 15000
$ws.Cells.Item(62, 11).Value = 15000
$ws.Cells.Item(62, 13).Value = -14314
$ws.Cells.Item(65, 8).Value = 27000
$ws.Cells.Item(65, 9).Value = 15000
$ws.Cells.Item(65, 11).Value = 45000
$ws.Cells.Item(65, 13).Value = -41568
$ws.Cells.Item(70, 9).Value = 6366.5
$ws.Cells.Item(70, 10).Value = 2000
$ws.Cells.Item(70, 11).Value = 6366.5
$ws.Cells.Item(70, 12).Value = 2000
$ws.Cells.Item(70, 13).Value = -6096.5
$ws.Cells.Item(70, 14).Value = -2540
$ws.Cells.Item(73, 9).Value = 6366.5
$ws.Cells.Item(73, 10).Value = 2000
$ws.Cells.Item(73, 11).Value = 6366.5
$ws.Cells.Item(73, 12).Value = 2000
$ws.Cells.Item(73, 13).Value = -5430.5
$ws.Cells.Item(73, 14).Value = -3872
$ws.Cells.Item(80, 8).Value = 4739.7144
$ws.Cells.Item(80, 9).Value = 3300
$ws.Cells.Item(80, 11).Value = 3300
$ws.Cells.Item(80, 13).Value = -2302
$ws.Cells.Item(83, 8).Value = 4739.7144
$ws.Cells.Item(83, 9).Value = 3300
$ws.Cells.Item(83, 11).Value = 16500
$ws.Cells.Item(83, 13).Value = -11508
$ws.Cells.Item(122, 8).Value = 10442341
$ws.Cells.Item(122, 9).Value = 12527670
$ws.Cells.Item(122, 11).Value = 37583010
$ws.Cells.Item(122, 13).Value = -37580560
$ws.Cells.Item(132, 8).Value = 1178
$ws.Cells.Item(132, 9).Value = 1178
$ws.Cells.Item(132, 11).Value = 3534
$ws.Cells.Item(132, 13).Value = -1004

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 9655.571
$ws.Cells.Item(40, 10).Value = 9930
$ws.Cells.Item(40, 12).Value = 9930
$ws.Cells.Item(40, 14).Value = -10202

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 14).Value = $null
$ws.Cells.Item(81, 8).Value = 343.70587
$ws.Cells.Item(81, 9).Value = 353.07144
$ws.Cells.Item(81, 10).Value = 300
$ws.Cells.Item(81, 11).Value = 706.14288
$ws.Cells.Item(81, 12).Value = 600
$ws.Cells.Item(81, 13).Value = 354.85712
$ws.Cells.Item(81, 14).Value = -2722
$ws.Cells.Item(84, 8).Value = 343.70587
$ws.Cells.Item(84, 9).Value = 353.07144
$ws.Cells.Item(84, 10).Value = 300
$ws.Cells.Item(84, 11).Value = 3530.7144
$ws.Cells.Item(84, 12).Value = 3000
$ws.Cells.Item(84, 13).Value = 1773.2856
$ws.Cells.Item(84, 14).Value = -13608
$ws.Cells.Item(132, 8).Value = 747.6
$ws.Cells.Item(132, 9).Value = 539.4286
$ws.Cells.Item(132, 11).Value = 1618.2858
$ws.Cells.Item(132, 13).Value = 911.7142000000001
